$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.020.21'
$ws.Range("E2").Value = '  +6.43%  '

$ws.Range("D3").Value = '1.883.82'
$ws.Range("E3").Value = '  +5.53%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.07'
$ws.Range("E5").Value = '  +1.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9983'
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5004'
$ws.Range("E7").Value = '  +1.76%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '45.76'
$ws.Range("E8").Value = '  +7.90%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2862'
$ws.Range("E9").Value = '  +6.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06559'
$ws.Range("E10").Value = '  +4.53%  '

$ws.Range("D11").Value = '1.878.25'
$ws.Range("E11").Value = '  +5.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '17.21'
$ws.Range("E12").Value = '  +4.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07242'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6697'
$ws.Range("E14").Value = '  +6.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.05'
$ws.Range("E15").Value = '  +6.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.823'
$ws.Range("E16").Value = '  +3.54%  '

$ws.Range("D17").Value = '30.016.29'
$ws.Range("E17").Value = '  +6.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9989'
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.91'
$ws.Range("E19").Value = '  +7.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007537'
$ws.Range("E20").Value = '  +4.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9984'
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").Value = '2.121.29'
$ws.Range("E22").Value = '  +5.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.776'
$ws.Range("E23").Value = '  +4.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.531'
$ws.Range("E24").Value = '  +5.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.028'
$ws.Range("E25").Value = '  +3.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.61'
$ws.Range("E26").Value = '  +2.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '134.81'
$ws.Range("E27").Value = '  +22.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.75'
$ws.Range("E28").Value = '  +6.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.952'
$ws.Range("E29").Value = '  +4.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.370'
$ws.Range("E30").Value = '  -1.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.197'
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08657'
$ws.Range("E32").Value = '  +4.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.918'
$ws.Range("E33").Value = '  +4.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05047'
$ws.Range("E34").Value = '  +3.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.139'
$ws.Range("E35").Value = '  +5.45%  '

$ws.Range("E36").Value = '  +6.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.688'
$ws.Range("E37").Value = '  +2.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.295'
$ws.Range("E38").Value = '  +12.16%  '

$ws.Range("E39").Value = '  +4.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9645'
$ws.Range("E40").Value = '  +1.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01641'
$ws.Range("E41").Value = '  +5.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.084'
$ws.Range("E42").Value = '  +2.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.80'
$ws.Range("E43").Value = '  +4.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9990'
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4225'
$ws.Range("E45").Value = '  +6.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.443'
$ws.Range("E46").Value = '  +3.53%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1258'
$ws.Range("E47").Value = '  +3.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05645'
$ws.Range("E48").Value = '  +3.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.46'
$ws.Range("E49").Value = '  +5.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.260'
$ws.Range("E50").Value = '  +3.37%  '

$ws.Range("E51").Value = '  +6.89%  '

